# ============================================================================
# Edit: Add a new "Player Info" worksheet (as the first sheet) containing
# basic player details, and update the "ODI Batting" / "ODI Bowling" sheets
# so the MATCH_CARD_LINK column (full howstat.com URL) becomes a MATCH_CODE
# column (just the numeric match code extracted from the URL).
# ============================================================================

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------------
# 1. Insert the new "Player Info" worksheet before the current first sheet
#    ("ODI Batting"), so the final sheet order is:
#       Player Info, ODI Batting, ODI Bowling
#    NOTE: worksheet references here are resolved positionally, so any
#    variable captured *before* inserting/moving sheets can silently start
#    pointing at a different sheet afterwards. Re-fetch sheets by name
#    once the insertion below has happened.
# ----------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# Header row
$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108   # xlCenter
$headerRange.VerticalAlignment = -4160     # xlTop
$headerRange.Borders.LineStyle = 1

# Data row - keep values as plain text (matching the rest of the workbook,
# which stores every cell as inline/shared string text).
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "6619"
$playerInfo.Range("A2").ClearFormats()

$playerInfo.Range("B2").Value = "Keacy Uydess Carty"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Medium"

$playerInfo.Range("A1").Select()

# ----------------------------------------------------------------------------
# 2. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE, and replace each
#    URL in column D with just the trailing numeric MatchCode.
# ----------------------------------------------------------------------------
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2  = "4577"
    3  = "4580"
    4  = "4583"
    5  = "4592"
    6  = "4616"
    7  = "4624"
    8  = "4636"
    9  = "4639"
    10 = "4642"
}

foreach ($row in $battingCodes.Keys) {
    $cell = $battingSheet.Range("D" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$row]
    $cell.ClearFormats()
}

# ----------------------------------------------------------------------------
# 3. "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE, and replace each
#    URL in column B with just the trailing numeric MatchCode.
# ----------------------------------------------------------------------------
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "4583"
    3 = "4616"
}

foreach ($row in $bowlingCodes.Keys) {
    $cell = $bowlingSheet.Range("B" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$row]
    $cell.ClearFormats()
}

Write-Host "Edit applied."
